# actualizacion enero 15 2025 - 8:50
$wb = $excel.ActiveWorkbook

# --- Update payment data on sheet "ene2025" (pago1 / pago2 columns) ---
$ws = $wb.Worksheets.Item("ene2025")

$ws.Range("C4").Value = 65000
$ws.Range("C9").Value = 65000
$ws.Range("C10").Value = 65000
$ws.Range("D10").Value = 65000
$ws.Range("C15").Value = 65000
$ws.Range("C16").Value = 65000
$ws.Range("C17").Value = 65000
$ws.Range("C19").Value = 65000
$ws.Range("C20").Value = 65000
$ws.Range("C23").Value = 65000

# --- Update selection on "cuota-extra" sheet (no longer the active tab) ---
$wsExtra = $wb.Worksheets.Item("cuota-extra")
$wsExtra.Range("G27").Select()

# --- Make "ene2025" the active sheet/tab and set its selection ---
$ws.Activate()
$ws.Range("C23").Select()
